$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap owner names (column E) between row 3 and row 4
$e3 = $ws.Range("E3").Value2
$e4 = $ws.Range("E4").Value2
$ws.Range("E3").Value = $e4
$ws.Range("E4").Value = $e3

# Update selection to K11
$ws.Range("K11").Select()
